$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder match rows: columns F..V are shuffled among rows; A..E (Indice/pais/torneio/temporada/data_partida) stay put ---
$row31 = @("Spezia", 0, "Como", 1, 1.59, "30/08/2023 19:42", 2.11, "03/09/2023 18:29", 4.15, "30/08/2023 19:42", 3.26, "03/09/2023 18:29", 6, "30/08/2023 19:42", 4.04, "03/09/2023 18:29", "https://www.betexplorer.com/football/italy/serie-b/spezia-como/hbBIlI94/")
for ($i = 0; $i -lt $row31.Length; $i++) { $ws.Cells.Item(31, 6 + $i).Value = $row31[$i] }

$row32 = @("Lecco", 3, "Catanzaro", 4, 2.51, "02/09/2023 01:42", 4.04, "03/09/2023 18:29", 3.26, "02/09/2023 01:42", 3.32, "03/09/2023 18:29", 3.07, "02/09/2023 01:42", 2.08, "03/09/2023 18:29", "https://www.betexplorer.com/football/italy/serie-b/lecco-catanzaro/4xVNSzmo/")
for ($i = 0; $i -lt $row32.Length; $i++) { $ws.Cells.Item(32, 6 + $i).Value = $row32[$i] }

$row33 = @("Cremonese", 1, "Sampdoria", 1, 2.05, "30/08/2023 19:42", 2.05, "03/09/2023 18:29", 3.47, "30/08/2023 19:42", 3.35, "03/09/2023 18:29", 3.79, "30/08/2023 19:42", 4.13, "03/09/2023 18:29", "https://www.betexplorer.com/football/italy/serie-b/cremonese-sampdoria/UqvdXGvU/")
for ($i = 0; $i -lt $row33.Length; $i++) { $ws.Cells.Item(33, 6 + $i).Value = $row33[$i] }

$row43 = @("Catanzaro", 0, "Parma", 5, 2.78, "06/09/2023 01:13", 2.44, "17/09/2023 16:06", 3.1, "06/09/2023 01:13", 3.3, "17/09/2023 16:06", 2.86, "06/09/2023 01:13", 3.18, "17/09/2023 16:06", "https://www.betexplorer.com/football/italy/serie-b/catanzaro-parma/Mw3Zpz9T/")
for ($i = 0; $i -lt $row43.Length; $i++) { $ws.Cells.Item(43, 6 + $i).Value = $row43[$i] }

$row44 = @("Como", 2, "Ternana", 1, 1.97, "04/09/2023 18:12", 2.19, "17/09/2023 16:11", 3.65, "04/09/2023 18:12", 3.28, "17/09/2023 16:11", 3.98, "04/09/2023 18:12", 3.75, "17/09/2023 16:11", "https://www.betexplorer.com/football/italy/serie-b/como-ternana/CIMgtfHp/")
for ($i = 0; $i -lt $row44.Length; $i++) { $ws.Cells.Item(44, 6 + $i).Value = $row44[$i] }

$row57 = @("Venezia", 1, "Palermo", 3, 2.24, "23/09/2023 13:13", 2.49, "26/09/2023 20:28", 3.46, "23/09/2023 13:13", 3.25, "26/09/2023 20:28", 3.37, "23/09/2023 13:13", 3.14, "26/09/2023 20:28", "https://www.betexplorer.com/football/italy/serie-b/venezia-palermo/ET3DbEfQ/")
for ($i = 0; $i -lt $row57.Length; $i++) { $ws.Cells.Item(57, 6 + $i).Value = $row57[$i] }

$row58 = @("Sudtirol", 0, "Modena", 0, 2.79, "23/09/2023 13:13", 2.86, "26/09/2023 20:29", 3.1, "23/09/2023 13:13", 2.8, "26/09/2023 20:29", 2.8, "23/09/2023 13:13", 3.11, "26/09/2023 20:29", "https://www.betexplorer.com/football/italy/serie-b/sudtirol-modena/dh98afuK/")
for ($i = 0; $i -lt $row58.Length; $i++) { $ws.Cells.Item(58, 6 + $i).Value = $row58[$i] }

$row59 = @("Spezia", 0, "Brescia", 0, 1.86, "25/09/2023 01:42", 1.83, "26/09/2023 20:29", 3.65, "25/09/2023 01:42", 3.58, "26/09/2023 20:29", 4.51, "25/09/2023 01:42", 4.89, "26/09/2023 20:29", "https://www.betexplorer.com/football/italy/serie-b/spezia-brescia/dUUCKR0a/")
for ($i = 0; $i -lt $row59.Length; $i++) { $ws.Cells.Item(59, 6 + $i).Value = $row59[$i] }

$row60 = @("Reggiana", 0, "Pisa", 0, 2.76, "25/09/2023 01:12", 2.33, "26/09/2023 20:29", 3.1, "25/09/2023 01:12", 3.08, "26/09/2023 20:29", 2.88, "25/09/2023 01:12", 3.64, "26/09/2023 20:29", "https://www.betexplorer.com/football/italy/serie-b/reggiana-pisa/I1A40zQD/")
for ($i = 0; $i -lt $row60.Length; $i++) { $ws.Cells.Item(60, 6 + $i).Value = $row60[$i] }

$row65 = @("Modena", 1, "Venezia", 3, 2.44, "26/09/2023 19:43", 2.47, "30/09/2023 13:26", 3.34, "26/09/2023 19:43", 3.25, "30/09/2023 13:26", 3.03, "26/09/2023 19:43", 3.17, "30/09/2023 13:26", "https://www.betexplorer.com/football/italy/serie-b/modena-venezia/Uadejj26/")
for ($i = 0; $i -lt $row65.Length; $i++) { $ws.Cells.Item(65, 6 + $i).Value = $row65[$i] }

$row66 = @("FeralpiSalo", 1, "Spezia", 2, 3.14, "28/09/2023 03:42", 3.67, "30/09/2023 13:57", 3.17, "28/09/2023 03:42", 3.18, "30/09/2023 13:46", 2.48, "28/09/2023 03:42", 2.27, "30/09/2023 13:46", "https://www.betexplorer.com/football/italy/serie-b/feralpisalo-spezia/j9ciiWn0/")
for ($i = 0; $i -lt $row66.Length; $i++) { $ws.Cells.Item(66, 6 + $i).Value = $row66[$i] }

$row67 = @("Brescia", 1, "Ascoli", 1, 2.15, "26/09/2023 19:43", 2.35, "30/09/2023 13:58", 3.31, "26/09/2023 19:43", 3.09, "30/09/2023 13:53", 3.67, "26/09/2023 19:43", 3.59, "30/09/2023 13:58", "https://www.betexplorer.com/football/italy/serie-b/brescia-ascoli/zLTGJoG5/")
for ($i = 0; $i -lt $row67.Length; $i++) { $ws.Cells.Item(67, 6 + $i).Value = $row67[$i] }

$row68 = @("Pisa", 1, "Cosenza", 2, 1.82, "26/09/2023 19:43", 2.04, "30/09/2023 13:57", 3.74, "26/09/2023 19:43", 3.24, "30/09/2023 13:57", 4.46, "26/09/2023 19:43", 4.33, "30/09/2023 13:57", "https://www.betexplorer.com/football/italy/serie-b/pisa-cosenza/xQa3lUWI/")
for ($i = 0; $i -lt $row68.Length; $i++) { $ws.Cells.Item(68, 6 + $i).Value = $row68[$i] }

$row70 = @("Sampdoria", 1, "Catanzaro", 2, 2.05, "28/09/2023 03:42", 2.26, "01/10/2023 16:12", 3.6, "28/09/2023 03:42", 3.31, "01/10/2023 16:11", 3.64, "28/09/2023 03:42", 3.52, "01/10/2023 16:12", "https://www.betexplorer.com/football/italy/serie-b/sampdoria-catanzaro/vVwnDD9J/")
for ($i = 0; $i -lt $row70.Length; $i++) { $ws.Cells.Item(70, 6 + $i).Value = $row70[$i] }

$row73 = @("Palermo", 2, "Sudtirol", 1, 1.74, "26/09/2023 19:43", 1.79, "01/10/2023 16:07", 3.72, "26/09/2023 19:43", 3.48, "01/10/2023 16:14", 5.03, "26/09/2023 19:43", 5.48, "01/10/2023 16:14", "https://www.betexplorer.com/football/italy/serie-b/palermo-sudtirol/bZ0akAHC/")
for ($i = 0; $i -lt $row73.Length; $i++) { $ws.Cells.Item(73, 6 + $i).Value = $row73[$i] }

$row77 = @("Reggiana", 1, "Bari", 1, 2.57, "02/10/2023 06:12", 2.81, "07/10/2023 13:57", 3.21, "02/10/2023 06:12", 3.09, "07/10/2023 13:57", 3.03, "02/10/2023 06:12", 2.87, "07/10/2023 13:57", "https://www.betexplorer.com/football/italy/serie-b/reggiana-bari/GWnV2ToC/")
for ($i = 0; $i -lt $row77.Length; $i++) { $ws.Cells.Item(77, 6 + $i).Value = $row77[$i] }

$row78 = @("Cosenza", 3, "Lecco", 0, 1.74, "02/10/2023 06:12", 1.89, "07/10/2023 13:33", 3.88, "02/10/2023 06:12", 3.48, "07/10/2023 13:34", 4.82, "02/10/2023 06:12", 4.69, "07/10/2023 12:45", "https://www.betexplorer.com/football/italy/serie-b/cosenza-lecco/86qN4kHa/")
for ($i = 0; $i -lt $row78.Length; $i++) { $ws.Cells.Item(78, 6 + $i).Value = $row78[$i] }

$row79 = @("Modena", 0, "Palermo", 2, 2.29, "01/10/2023 15:42", 2.9, "07/10/2023 13:52", 3.41, "01/10/2023 15:42", 3.12, "07/10/2023 13:57", 3.3, "01/10/2023 15:42", 2.77, "07/10/2023 13:56", "https://www.betexplorer.com/football/italy/serie-b/modena-palermo/2TrR39W5/")
for ($i = 0; $i -lt $row79.Length; $i++) { $ws.Cells.Item(79, 6 + $i).Value = $row79[$i] }

$row80 = @("Ascoli", 1, "Sampdoria", 1, 2.55, "02/10/2023 06:12", 2.69, "07/10/2023 16:12", 3.25, "02/10/2023 06:12", 3.1, "07/10/2023 16:12", 3.01, "02/10/2023 06:12", 3, "07/10/2023 16:12", "https://www.betexplorer.com/football/italy/serie-b/ascoli-sampdoria/ADUH7iXt/")
for ($i = 0; $i -lt $row80.Length; $i++) { $ws.Cells.Item(80, 6 + $i).Value = $row80[$i] }

$row81 = @("Cittadella", 2, "Ternana", 2, 2.05, "01/10/2023 15:42", 2.33, "07/10/2023 16:11", 3.46, "01/10/2023 15:42", 3.06, "07/10/2023 16:13", 3.92, "01/10/2023 15:42", 3.67, "07/10/2023 16:11", "https://www.betexplorer.com/football/italy/serie-b/cittadella-ternana/l4TL6Bnm/")
for ($i = 0; $i -lt $row81.Length; $i++) { $ws.Cells.Item(81, 6 + $i).Value = $row81[$i] }

$row82 = @("Venezia", 3, "Parma", 2, 2.35, "01/10/2023 17:43", 2.76, "07/10/2023 16:13", 3.38, "01/10/2023 17:43", 3.22, "07/10/2023 16:13", 3.14, "01/10/2023 17:43", 2.82, "07/10/2023 16:13", "https://www.betexplorer.com/football/italy/serie-b/venezia-parma/Q1ULPjgP/")
for ($i = 0; $i -lt $row82.Length; $i++) { $ws.Cells.Item(82, 6 + $i).Value = $row82[$i] }

$row86 = @("Ternana", 0, "Brescia", 1, 2.08, "09/10/2023 16:12", 2.01, "21/10/2023 13:59", 3.38, "09/10/2023 16:12", 3.33, "21/10/2023 13:59", 3.92, "09/10/2023 16:12", 4.33, "21/10/2023 13:59", "https://www.betexplorer.com/football/italy/serie-b/ternana-brescia/WnNPHPpI/")
for ($i = 0; $i -lt $row86.Length; $i++) { $ws.Cells.Item(86, 6 + $i).Value = $row86[$i] }

$row87 = @("Pisa", 2, "Cittadella", 1, 2.17, "09/10/2023 16:12", 2.3, "21/10/2023 13:55", 3.2, "09/10/2023 16:12", 2.91, "21/10/2023 13:53", 3.88, "09/10/2023 16:12", 4.02, "21/10/2023 13:55", "https://www.betexplorer.com/football/italy/serie-b/pisa-cittadella/jk1DF6VH/")
for ($i = 0; $i -lt $row87.Length; $i++) { $ws.Cells.Item(87, 6 + $i).Value = $row87[$i] }

$row89 = @("Cremonese", 0, "Sudtirol", 1, 1.75, "10/10/2023 01:12", 1.59, "21/10/2023 13:55", 3.79, "10/10/2023 01:12", 3.97, "21/10/2023 13:55", 5.01, "10/10/2023 01:12", 6.62, "21/10/2023 13:55", "https://www.betexplorer.com/football/italy/serie-b/cremonese-sudtirol/M99dJlWh/")
for ($i = 0; $i -lt $row89.Length; $i++) { $ws.Cells.Item(89, 6 + $i).Value = $row89[$i] }

$row90 = @("Bari", 1, "Modena", 1, 2.19, "09/10/2023 16:12", 2.43, "21/10/2023 13:58", 3.25, "09/10/2023 16:12", 2.97, "21/10/2023 13:57", 3.74, "09/10/2023 16:12", 3.58, "21/10/2023 13:58", "https://www.betexplorer.com/football/italy/serie-b/bari-modena/SAGmLA1t/")
for ($i = 0; $i -lt $row90.Length; $i++) { $ws.Cells.Item(90, 6 + $i).Value = $row90[$i] }

$row103 = @("Venezia", 2, "Pisa", 1, 2.04, "24/10/2023 18:42", 2.15, "29/10/2023 16:14", 3.45, "24/10/2023 18:42", 3.27, "29/10/2023 16:14", 3.85, "24/10/2023 18:42", 3.9, "29/10/2023 16:14", "https://www.betexplorer.com/football/italy/serie-b/venezia-pisa/KpgOAHqj/")
for ($i = 0; $i -lt $row103.Length; $i++) { $ws.Cells.Item(103, 6 + $i).Value = $row103[$i] }

$row106 = @("Brescia", 1, "Bari", 2, 2.65, "24/10/2023 18:42", 2.61, "29/10/2023 16:14", 3.05, "24/10/2023 18:42", 2.85, "29/10/2023 16:12", 3.01, "24/10/2023 18:42", 3.4, "29/10/2023 16:14", "https://www.betexplorer.com/football/italy/serie-b/brescia-bari/Orog0cFp/")
for ($i = 0; $i -lt $row106.Length; $i++) { $ws.Cells.Item(106, 6 + $i).Value = $row106[$i] }

$row107 = @("Bari", 1, "Ascoli", 0, 1.98, "29/10/2023 16:42", 2.14, "04/11/2023 13:59", 3.35, "29/10/2023 16:42", 3.02, "04/11/2023 13:59", 4.34, "29/10/2023 16:42", 4.34, "04/11/2023 13:59", "https://www.betexplorer.com/football/italy/serie-b/bari-ascoli/byhS9yad/")
for ($i = 0; $i -lt $row107.Length; $i++) { $ws.Cells.Item(107, 6 + $i).Value = $row107[$i] }

$row108 = @("Catanzaro", 1, "Modena", 2, 2.05, "29/10/2023 16:42", 2.37, "04/11/2023 13:58", 3.49, "29/10/2023 16:42", 3.32, "04/11/2023 13:58", 3.89, "29/10/2023 16:42", 3.27, "04/11/2023 13:58", "https://www.betexplorer.com/football/italy/serie-b/catanzaro-modena/CvtAYKyc/")
for ($i = 0; $i -lt $row108.Length; $i++) { $ws.Cells.Item(108, 6 + $i).Value = $row108[$i] }

$row109 = @("Cittadella", 3, "Brescia", 2, 1.98, "29/10/2023 16:42", 2.12, "04/11/2023 13:55", 3.34, "29/10/2023 16:42", 3.01, "04/11/2023 13:52", 4.36, "29/10/2023 16:42", 4.47, "04/11/2023 13:55", "https://www.betexplorer.com/football/italy/serie-b/cittadella-brescia/zJlcaHUj/")
for ($i = 0; $i -lt $row109.Length; $i++) { $ws.Cells.Item(109, 6 + $i).Value = $row109[$i] }

$row110 = @("Cosenza", 1, "FeralpiSalo", 1, 2.01, "28/10/2023 15:13", 1.92, "04/11/2023 13:57", 3.43, "28/10/2023 15:13", 3.37, "04/11/2023 13:57", 4.09, "28/10/2023 15:13", 4.73, "04/11/2023 13:57", "https://www.betexplorer.com/football/italy/serie-b/cosenza-feralpisalo/OEqEXvi3/")
for ($i = 0; $i -lt $row110.Length; $i++) { $ws.Cells.Item(110, 6 + $i).Value = $row110[$i] }

$row111 = @("Pisa", 1, "Como", 1, 2.28, "29/10/2023 16:42", 2.66, "04/11/2023 13:57", 3.45, "29/10/2023 16:42", 2.93, "04/11/2023 13:58", 3.23, "29/10/2023 16:42", 3.23, "04/11/2023 13:57", "https://www.betexplorer.com/football/italy/serie-b/pisa-como/E3xRUxyM/")
for ($i = 0; $i -lt $row111.Length; $i++) { $ws.Cells.Item(111, 6 + $i).Value = $row111[$i] }

$row112 = @("Ternana", 0, "Venezia", 1, 2.49, "29/10/2023 16:42", 2.82, "04/11/2023 13:59", 3.39, "29/10/2023 16:42", 3.26, "04/11/2023 13:59", 2.98, "29/10/2023 16:42", 2.73, "04/11/2023 13:59", "https://www.betexplorer.com/football/italy/serie-b/ternana-venezia/CvI1OzEk/")
for ($i = 0; $i -lt $row112.Length; $i++) { $ws.Cells.Item(112, 6 + $i).Value = $row112[$i] }

# --- Append new rows 114-116 (copy formatting from row 113, then set values) ---
$ws.Range("A113:V113").Copy($ws.Range("A114:V114"))
$ws.Cells.Item(114, 1).Value = 113
$ws.Cells.Item(114, 2).Value = "italy"
$ws.Cells.Item(114, 3).Value = "serie-b"
$ws.Cells.Item(114, 4).Value = "2023-2024"
$ws.Cells.Item(114, 5).Value = 45235.67708333334
$ws.Cells.Item(114, 6).Value = "Parma"
$ws.Cells.Item(114, 7).Value = 2
$ws.Cells.Item(114, 8).Value = "Sudtirol"
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 1.61
$ws.Cells.Item(114, 11).Value = "29/10/2023 16:43"
$ws.Cells.Item(114, 12).Value = 1.72
$ws.Cells.Item(114, 13).Value = "05/11/2023 16:06"
$ws.Cells.Item(114, 14).Value = 3.98
$ws.Cells.Item(114, 15).Value = "29/10/2023 16:43"
$ws.Cells.Item(114, 16).Value = 3.65
$ws.Cells.Item(114, 17).Value = "05/11/2023 16:06"
$ws.Cells.Item(114, 18).Value = 6.1
$ws.Cells.Item(114, 19).Value = "29/10/2023 16:43"
$ws.Cells.Item(114, 20).Value = 5.67
$ws.Cells.Item(114, 21).Value = "05/11/2023 16:06"
$ws.Cells.Item(114, 22).Value = "https://www.betexplorer.com/football/italy/serie-b/parma-sudtirol/pUXMVILF/"

$ws.Range("A113:V113").Copy($ws.Range("A115:V115"))
$ws.Cells.Item(115, 1).Value = 114
$ws.Cells.Item(115, 2).Value = "italy"
$ws.Cells.Item(115, 3).Value = "serie-b"
$ws.Cells.Item(115, 4).Value = "2023-2024"
$ws.Cells.Item(115, 5).Value = 45235.67708333334
$ws.Cells.Item(115, 6).Value = "Cremonese"
$ws.Cells.Item(115, 7).Value = 3
$ws.Cells.Item(115, 8).Value = "Spezia"
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 10).Value = 1.78
$ws.Cells.Item(115, 11).Value = "29/10/2023 16:43"
$ws.Cells.Item(115, 12).Value = 1.99
$ws.Cells.Item(115, 13).Value = "05/11/2023 16:06"
$ws.Cells.Item(115, 14).Value = 3.78
$ws.Cells.Item(115, 15).Value = "29/10/2023 16:43"
$ws.Cells.Item(115, 16).Value = 3.36
$ws.Cells.Item(115, 17).Value = "05/11/2023 16:08"
$ws.Cells.Item(115, 18).Value = 4.86
$ws.Cells.Item(115, 19).Value = "29/10/2023 16:43"
$ws.Cells.Item(115, 20).Value = 4.36
$ws.Cells.Item(115, 21).Value = "05/11/2023 16:08"
$ws.Cells.Item(115, 22).Value = "https://www.betexplorer.com/football/italy/serie-b/cremonese-spezia/6wZIWb69/"

$ws.Range("A113:V113").Copy($ws.Range("A116:V116"))
$ws.Cells.Item(116, 1).Value = 115
$ws.Cells.Item(116, 2).Value = "italy"
$ws.Cells.Item(116, 3).Value = "serie-b"
$ws.Cells.Item(116, 4).Value = "2023-2024"
$ws.Cells.Item(116, 5).Value = 45235.67708333334
$ws.Cells.Item(116, 6).Value = "Reggiana"
$ws.Cells.Item(116, 7).Value = 1
$ws.Cells.Item(116, 8).Value = "Lecco"
$ws.Cells.Item(116, 9).Value = 1
$ws.Cells.Item(116, 10).Value = 1.77
$ws.Cells.Item(116, 11).Value = "30/10/2023 07:12"
$ws.Cells.Item(116, 12).Value = 1.74
$ws.Cells.Item(116, 13).Value = "05/11/2023 13:45"
$ws.Cells.Item(116, 14).Value = 3.72
$ws.Cells.Item(116, 15).Value = "30/10/2023 07:12"
$ws.Cells.Item(116, 16).Value = 3.81
$ws.Cells.Item(116, 17).Value = "05/11/2023 15:02"
$ws.Cells.Item(116, 18).Value = 4.97
$ws.Cells.Item(116, 19).Value = "30/10/2023 07:12"
$ws.Cells.Item(116, 20).Value = 5.09
$ws.Cells.Item(116, 21).Value = "05/11/2023 16:07"
$ws.Cells.Item(116, 22).Value = "https://www.betexplorer.com/football/italy/serie-b/reggiana-lecco/hfyVTdjS/"

$ws.Range("A1").Select() | Out-Null
